$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.01%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.79%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.54%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06217"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.32%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.754"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.40%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8514"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.34%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9145"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.40%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.19%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04906"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.44%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07091"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.07%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.50%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09047"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.35%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001532"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.50%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006166"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.62%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006100"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.12%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.16%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.180"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.92%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.165"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.47%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.10%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.083"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.11%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04235"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.61%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.24%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004078"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.27%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.04%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.40%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03958"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.19%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1114"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.14%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004132"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.14%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-15.11%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.33%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.04%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2483"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "87.96%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.04%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
